$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SS1_Pu_GHS")

$cols = @("AL","AR","AS","AT","AU","AV","AW","AX","AY","AZ","BB","BC","BE","BF","BH")

for ($row = 15; $row -le 27; $row++) {
    foreach ($col in $cols) {
        $addr = "$col$row"
        $ws.Range($addr).ClearContents()
    }
}
